# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect freshly-generated stats (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1070
    5  = 81
    7  = 53
    8  = 11088
    9  = 4265
    10 = 24
    11 = 19
    13 = 2491
    14 = 1066
    15 = 79
    16 = 12
    17 = 152
    18 = 476
    19 = 11200
    20 = 11033
    21 = 14
    24 = 10
    25 = 33
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
